$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab_4a_Indikatorenblätter")

# Row 29 - Anteil erneuerbarer Energien am Brutto-Endenergieverbrauch
$ws.Range("D29").Value = "Anteil erneuerbarer Energien am Bruttoendenergieverbrauch"
$ws.Range("F29").Value = "Der Indikator setzt die Erzeugung erneuerbarer Energien in Relation zum Bruttoendenergieverbrauch. Der Brutto-Endenergieverbrauch umfasst den Energieverbrauch beim Letztverbraucher, die Übertragungsverluste sowie den Eigenverbrauch der Energiegewinnungsbereiche."

# Row 34 - Verhältnis der Bruttoanlageinvestitionen zum BIP
$ws.Range("D34").Value = "Verhältnis der Bruttoanlageinvestitionen zum Bruttoinlandsprodukt (Investitionsquote)"

# Row 38 - Freiwillige Nachhaltigkeitsberichterstattung
$ws.Range("D38").Value = "Freiwillige Nachhaltigkeitsberichterstattung von Unternehmen nach dem deutschen Nachhaltigkeitskodex (DNK)"
$ws.Range("E38").Value = "XXXFreiwillige Nachhaltigkeitsberichterstattung von Unternehmen nach dem deutschen Nachhaltigkeitskodex (DNK)"

# Row 40 - Breitbandausbau
$ws.Range("D40").Value = "Breitbandausbau - Anteil der Haushalte mit Zugang zu Gigabit-Breitbandversorgung"
$ws.Range("E40").Value = "XXXRollout of broadband"

# Row 48 - Überlastung durch Wohnkosten
$ws.Range("D48").Value = "Anteil der Personen mit hohen Wohnkosten"
$ws.Range("E48").Value = "XXXHousing cost overload"

# Row 52 - Umweltmanagement EMAS
$ws.Range("D52").Value = "Umweltmanagementsystem EMAS"
$ws.Range("E52").Value = "EMAS Eco-Management and Audit Scheme"

# Row 58 - Stickstoffeintrag über die Zuflüsse in die Nord- und Ostsee
$ws.Range("D58").Value = "Nährstoffeinträge in Küsten- und Meeresgewässer - Stickstoffeintrag über die Zuflüsse in die Ost- und Nordsee"
$ws.Range("E58").Value = "XXXNährstoffeinträge in Küsten- und Meeresgewässer - Stickstoffeintrag über die Zuflüsse in die Ost- und Nordsee"

# Row 63 - Bilaterale Beiträge der deutschen internationalen Kooperation ...
$ws.Range("D63").Value = "Bilaterale Beiträge der deutschen internationalen Kooperation zum Schutz, nachhaltiger Nutzung und Wiederherstellung von Land"
$ws.Range("E63").Value = "XXXBilaterale Beiträge der deutschen internationalen Kooperation zum Schutz, nachhaltiger Nutzung und Wiederherstellung von Land"
